# Updated Lookup tables with symbol changes
# Remove the SPICEJET row (row 44, which holds "SPICEJET" in both columns A and B)
# and shift everything below it up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $null
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "SPICEJET") {
        $target = $r
        break
    }
}

if ($target -ne $null) {
    $ws.Rows.Item($target).Delete()
}

# Update the view to match the new cursor position
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("H37").Select()
